{"js": "// Move the empty \"_GoBack\" bookmark from after the \"Exercises\" run\n// to the very start of that same (first) paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Remove the existing \"_GoBack\" bookmark (currently located after the\n// \"Exercises\" text run).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Re-insert it at the start of the first paragraph (before the run\n// containing \"Exercises\").\nconst startRange = firstParagraph.getRange(\"Start\");\nstartRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Move the empty \"_GoBack\" bookmark from after the \"Exercises\" run\n# to the very start of that same (first) paragraph.\n\n$d = $word.ActiveDocument\n\n# Remove the existing \"_GoBack\" bookmark (currently located after the\n# \"Exercises\" text run, at the end of the first paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Build a range that covers just the paragraph mark of the first\n# paragraph (i.e. its very last \"character\"). Anchoring the new\n# bookmark to the paragraph-mark range - rather than a plain\n# zero-length range at the paragraph's text start - reliably makes\n# Word store the bookmark immediately after the paragraph properties\n# and before the \"Exercises\" run, without splitting that run apart.\n$firstParagraph = $d.Paragraphs(1)\n$paragraphEnd = $firstParagraph.Range.End\n$markRange = $d.Range($paragraphEnd - 1, $paragraphEnd)\n\n$d.Bookmarks.Add(\"_GoBack\", $markRange)\n"}
